# Reverse the order of comma-separated names/emails in the "Recorded By"
# column (column G) for every data row on the active sheet.
#
# For example:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system"   -> "system, backup@backdoor.com, System"
#
# Cells that contain only a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value

    if ($text -notmatch ",") {
        continue
    }

    $parts = $text -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Manually reverse the token order (([array]::Reverse() is not
    # reliable in this runtime), by walking the list back-to-front.
    $reversed = @()
    for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
        $reversed += $trimmed[$i]
    }

    $newValue = [string]::Join(", ", $reversed)

    # Note: PowerShell's -eq/-ne operators on strings in this runtime are
    # case-insensitive, which would incorrectly skip cells that only
    # differ by letter case (e.g. "System" vs "system"). Always write the
    # recomputed value back; writing an unchanged value is harmless.
    $cell.Value2 = $newValue
}
